# Adds a 5th account column ("telefonia") to the CuentasBancarias sheet,
# matching the formatting already used by the existing "cuenta bancaria"
# and "cuenta contable" rows (copy the format from a sibling cell instead
# of inventing a brand-new font/style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "telefonia"
$ws.Range("E3").Value = "davivienda"
$ws.Range("E4").Value = "aho"
$ws.Range("E5").Value = "0010-1005-222"
$ws.Range("E6").Value = "0698.111006.1"

# Reuse existing cell formatting (number format + font) rather than create
# new style/font entries: E5 matches B5's style, E6 matches D6's style.
$ws.Range("B5").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("C9").Select()
